# Update the "Absent" column (H) values on the consolidated report sheet.
# H is derived from E (Real attendance): Absent = 1 when Real = 0, else 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
